{"js": "// Update the 25 \"NNN\u00d7N=\" multiplication prompts in the (single) table to\n// the new values, in document order. Several old prompts repeat (e.g.\n// \"892\u00d74=\" appears twice) but map to different new values depending on\n// position, so replacement must walk cells in order rather than doing a\n// global find/replace.\nconst replacements = [\n  \"173\u00d73=\", \"447\u00d74=\", \"691\u00d77=\", \"967\u00d78=\", \"747\u00d75=\",\n  \"485\u00d79=\", \"993\u00d79=\", \"998\u00d77=\", \"577\u00d72=\", \"689\u00d72=\",\n  \"949\u00d78=\", \"628\u00d73=\", \"987\u00d75=\", \"319\u00d76=\", \"369\u00d74=\",\n  \"863\u00d79=\", \"785\u00d77=\", \"429\u00d72=\", \"554\u00d79=\", \"824\u00d73=\",\n  \"968\u00d74=\", \"178\u00d79=\", \"210\u00d75=\", \"300\u00d78=\", \"377\u00d76=\",\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nlet cursor = 0;\n\nfor (let ti = 0; ti < tables.items.length && cursor < replacements.length; ti++) {\n  const table = tables.items[ti];\n  table.load(\"rowCount,values\");\n  await context.sync();\n\n  const colCount = table.values.length > 0 ? table.values[0].length : 0;\n\n  for (let r = 0; r < table.rowCount && cursor < replacements.length; r++) {\n    for (let c = 0; c < colCount && cursor < replacements.length; c++) {\n      const cellText = table.values[r][c];\n      if (cellText === \"\" || cellText === null || cellText === undefined) {\n        continue; // skip blank spacer cells\n      }\n\n      const cell = table.getCellOrNullObject(r, c);\n      cell.body.load(\"paragraphs/items\");\n      await context.sync();\n\n      if (cell.isNullObject) {\n        continue;\n      }\n\n      const para = cell.body.paragraphs.items[0];\n      const range = para.getRange();\n      range.insertText(replacements[cursor], Word.InsertLocation.replace);\n      cursor++;\n    }\n  }\n  await context.sync();\n}\n\nawait context.sync();\n", "ps1": "# Update the 25 \"NNN\u00d7N=\" multiplication prompts in the (single) table to\n# the new values, in document order. Several old prompts repeat (e.g.\n# \"892\u00d74=\" appears twice) but map to different new values depending on\n# position, so replacement walks cells in row-major order rather than\n# doing a single global find/replace.\n$replacements = @(\n    \"173\u00d73=\", \"447\u00d74=\", \"691\u00d77=\", \"967\u00d78=\", \"747\u00d75=\",\n    \"485\u00d79=\", \"993\u00d79=\", \"998\u00d77=\", \"577\u00d72=\", \"689\u00d72=\",\n    \"949\u00d78=\", \"628\u00d73=\", \"987\u00d75=\", \"319\u00d76=\", \"369\u00d74=\",\n    \"863\u00d79=\", \"785\u00d77=\", \"429\u00d72=\", \"554\u00d79=\", \"824\u00d73=\",\n    \"968\u00d74=\", \"178\u00d79=\", \"210\u00d75=\", \"300\u00d78=\", \"377\u00d76=\"\n)\n\n$d = $word.ActiveDocument\n$cursor = 0\n\nforeach ($t in $d.Tables) {\n    $rows = $t.Rows.Count\n    $cols = $t.Columns.Count\n\n    for ($r = 1; $r -le $rows; $r++) {\n        for ($c = 1; $c -le $cols; $c++) {\n            if ($cursor -ge $replacements.Count) { continue }\n\n            $cell = $t.Cell($r, $c)\n            $range = $cell.Range\n\n            # An empty cell's Range is just the end-of-cell mark (2 chars:\n            # paragraph mark + cell mark). Only touch cells that actually\n            # have visible text, and check the length BEFORE trimming --\n            # this engine's Range.Text on a collapsed (Start==End) range\n            # does not reliably reflect the post-trim content.\n            if ($range.Text.Length -gt 2) {\n                # Cell.Range includes the trailing cell-mark (and\n                # paragraph mark); trim it off so we only replace the\n                # visible text and keep the run/paragraph formatting\n                # intact.\n                $range.End = $range.End - 1\n                $range.Text = $replacements[$cursor]\n                $cursor++\n            }\n        }\n    }\n}\n"}
